$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: Price values in column D are plain text strings (e.g. "1.000", "31.300.18")
# that look numeric. Assigning them with a leading apostrophe forces Excel to keep
# them as literal text instead of auto-converting to a Double (which would strip
# trailing zeros / reformat using scientific notation).

$ws.Range("D2").Value = "'31.300.18"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").Value = "'1.966.94"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'248.37"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4901"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("D8").Value = "'44.76"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "'0.2990"
$ws.Range("E9").Value = "  +3.47%  "
$ws.Range("D10").Value = "'0.06870"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").Value = "'19.32"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "'107.25"
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07781"
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.940.16"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").Value = "'5.460"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "'0.7186"
$ws.Range("E16").Value = "  +7.13%  "
$ws.Range("D17").Value = "'288.92"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").Value = "'31.295.40"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "'0.000007805"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").Value = "'5.641"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").Value = "'2.197.56"
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'6.644"
$ws.Range("E25").Value = "  +3.33%  "
$ws.Range("D26").Value = "'10.04"
$ws.Range("E26").Value = "  +5.98%  "
$ws.Range("D27").Value = "'168.99"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").Value = "'20.09"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").Value = "'2.208"
$ws.Range("E29").Value = "  +5.80%  "
$ws.Range("D30").Value = "'0.1071"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "'1.448"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").Value = "'4.852"
$ws.Range("E32").Value = "  +19.70%  "
$ws.Range("D33").Value = "'4.530"
$ws.Range("E33").Value = "  +9.78%  "
$ws.Range("D34").Value = "'0.05074"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("D35").Value = "'0.7732"
$ws.Range("E35").Value = "  +5.31%  "
$ws.Range("D36").Value = "'1.177"
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("D37").Value = "'0.02062"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Value = "'2.731"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").Value = "'2.718"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("D40").Value = "'2.149"
$ws.Range("E40").Value = "  +6.33%  "
$ws.Range("D41").Value = "'6.428"
$ws.Range("E41").Value = "  +10.27%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8894"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4498"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'110.19"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'73.89"
$ws.Range("E45").Value = "  +6.58%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'7.551"
$ws.Range("E47").Value = "  +4.82%  "
$ws.Range("D48").Value = "'999.69"
$ws.Range("E48").Value = "  +18.76%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1277"
$ws.Range("E49").Value = "  +4.28%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.399"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").Value = "'36.25"
$ws.Range("E51").Value = "  +4.39%  "
